# Add data for 2022-02-17: rename sheet/title to "through 02-09" and
# bump the affected totals (new Feb-2022 entry + a correction to 2021's
# September figure).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet ("Through 2022-02-08" -> "Through 2022-02-09")
$ws.Name = "Through 2022-02-09"

# Update the column header label in A1's shared string ("2022 (through 02-08)" -> "2022 (through 02-09)")
$ws.Range("I1").Value = "2022 (through 02-09)"

# February / Total column (I3): 39 -> 40
$ws.Range("I3").Value = 40

# September / 2021 column (H10): 177 -> 178
$ws.Range("H10").Value = 178

# Total row / 2021 column (H14): 1852 -> 1853
$ws.Range("H14").Value = 1853

# Total row / Total column (I14): 201 -> 202
$ws.Range("I14").Value = 202
